# Add a new "2022-Q3" sheet (as the 2nd tab, right after "总计") with its
# own fund-holding table, and insert a corresponding summary row into the
# "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)          # "总计"
$q2Sheet    = $wb.Worksheets.Item(2)           # "2022-Q2" (template for layout/styles)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying the "2022-Q2" sheet so it
#    inherits the exact same layout / number formats / column widths,
#    then overwrite its data with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Force the fund-code / percentage-like text columns to stay text so
# leading zeros (e.g. "006679") and the textual number look of the
# source data survive the round trip.
$q3Sheet.Range("B2:B5").NumberFormat = "@"
$q3Sheet.Range("D2:G5").NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "006679"
$q3Sheet.Range("C2").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 A"
$q3Sheet.Range("D2").Value = "11.73"
$q3Sheet.Range("E2").Value = "93.96"
$q3Sheet.Range("F2").Value = "4.06"
$q3Sheet.Range("G2").Value = "0.4762"
$q3Sheet.Range("H2").Value = 7

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "162719"
$q3Sheet.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
$q3Sheet.Range("D3").Value = "11.73"
$q3Sheet.Range("E3").Value = "93.96"
$q3Sheet.Range("F3").Value = "4.06"
$q3Sheet.Range("G3").Value = "0.4762"
$q3Sheet.Range("H3").Value = 7

$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "006680"
$q3Sheet.Range("C4").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 C"
$q3Sheet.Range("D4").Value = "5.92"
$q3Sheet.Range("E4").Value = "93.96"
$q3Sheet.Range("F4").Value = "4.06"
$q3Sheet.Range("G4").Value = "0.2404"
$q3Sheet.Range("H4").Value = 7

$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "004243"
$q3Sheet.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
$q3Sheet.Range("D5").Value = "-11.74"
$q3Sheet.Range("E5").Value = "93.96"
$q3Sheet.Range("F5").Value = "4.06"
$q3Sheet.Range("G5").Value = "-0.4766"
$q3Sheet.Range("H5").Value = 7

# ---------------------------------------------------------------------
# 2) Insert a new summary row for 2022-Q3 into the "总计" sheet, pushing
#    the existing quarters down by one row.
# ---------------------------------------------------------------------
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Range("B$dest").Value = $totalSheet.Range("B$r").Value2
    $totalSheet.Range("C$dest").Value = $totalSheet.Range("C$r").Value2
    $totalSheet.Range("D$dest").Value = $totalSheet.Range("D$r").Value2
    $totalSheet.Range("A$dest").Value = $dest - 2
}

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.72

# Keep "总计" as the active tab, matching the original workbook view state.
$totalSheet.Activate()
